$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "95.752.05"
$ws.Range("E2").Value = "  -1.31%  "

Set-TextCell $ws.Range("D3") "3.612.91"
$ws.Range("E3").Value = "  -1.98%  "

Set-TextCell $ws.Range("D4") "2.73"
$ws.Range("E4").Value = "  +28.39%  "

$ws.Range("E5").Value = "  +0.12%  "

Set-TextCell $ws.Range("D6") "223.32"
$ws.Range("E6").Value = "  -5.08%  "

Set-TextCell $ws.Range("D7") "640.79"
$ws.Range("E7").Value = "  -2.14%  "

Set-TextCell $ws.Range("D8") "0.424"
$ws.Range("E8").Value = "  -3.04%  "

Set-TextCell $ws.Range("D9") "1.20"
$ws.Range("E9").Value = "  +9.88%  "

$ws.Range("E10").Value = "  +0.06%  "

Set-TextCell $ws.Range("D11") "3.610.10"
$ws.Range("E11").Value = "  -2.02%  "

Set-TextCell $ws.Range("D12") "48.51"
$ws.Range("E12").Value = "  +9.22%  "

$ws.Range("E13").Value = "  +4.56%  "

$ws.Range("E14").Value = "  -5.95%  "

Set-TextCell $ws.Range("D15") "6.53"
$ws.Range("E15").Value = "  -3.95%  "

Set-TextCell $ws.Range("D16") "4.287.94"
$ws.Range("E16").Value = "  -2.00%  "

Set-TextCell $ws.Range("D17") "24.77"
$ws.Range("E17").Value = "  +33.79%  "

Set-TextCell $ws.Range("D18") "95.733.35"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("E19").Value = "  +5.13%  "

Set-TextCell $ws.Range("D20") "13.82"
$ws.Range("E20").Value = "  +6.70%  "

Set-TextCell $ws.Range("D21") "3.617.80"
$ws.Range("E21").Value = "  -1.99%  "

Set-TextCell $ws.Range("D22") "0.292"
$ws.Range("E22").Value = "  +44.83%  "

Set-TextCell $ws.Range("D23") "0.537"
$ws.Range("E23").Value = "  +0.25%  "

Set-TextCell $ws.Range("D24") "137.03"
$ws.Range("E24").Value = "  +23.73%  "

Set-TextCell $ws.Range("D25") "526.00"
$ws.Range("E25").Value = "  +1.70%  "

Set-TextCell $ws.Range("D26") "3.26"
$ws.Range("E26").Value = "  -4.42%  "

Set-TextCell $ws.Range("D27") "0.0000202"
$ws.Range("E27").Value = "  -7.79%  "

Set-TextCell $ws.Range("D28") "6.85"
$ws.Range("E28").Value = "  -0.07%  "

Set-TextCell $ws.Range("D29") "3.787.21"
$ws.Range("E29").Value = "  -2.47%  "

Set-TextCell $ws.Range("D30") "12.93"
$ws.Range("E30").Value = "  -3.04%  "

Set-TextCell $ws.Range("D31") "13.27"
$ws.Range("E31").Value = "  +5.92%  "

Set-TextCell $ws.Range("D32") "3.14"
$ws.Range("E32").Value = "  +5.09%  "

$ws.Range("E33").Value = "  +0.12%  "

Set-TextCell $ws.Range("D34") "0.635"
$ws.Range("E34").Value = "  +7.84%  "

Set-TextCell $ws.Range("D35") "33.39"
$ws.Range("E35").Value = "  +2.61%  "

Set-TextCell $ws.Range("D36") "0.183"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("E37").Value = "  +1.42%  "

$ws.Range("E38").Value = "  +0.20%  "

Set-TextCell $ws.Range("D39") "0.536"
$ws.Range("E39").Value = "  +8.83%  "

$ws.Range("E40").Value = "  +0.00%  "

Set-TextCell $ws.Range("D41") "7.22"
$ws.Range("E41").Value = "  +6.23%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D42") "0.0538"
$ws.Range("E42").Value = "  +20.48%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws.Range("D43") "589.81"
$ws.Range("E43").Value = "  -6.27%  "

Set-TextCell $ws.Range("D44") "8.38"
$ws.Range("E44").Value = "  -3.68%  "

Set-TextCell $ws.Range("D45") "41.29"
$ws.Range("E45").Value = "  +3.25%  "

Set-TextCell $ws.Range("D46") "1.00"
$ws.Range("E46").Value = "  +5.51%  "

Set-TextCell $ws.Range("D47") "0.158"
$ws.Range("E47").Value = "  -4.69%  "

Set-TextCell $ws.Range("D48") "1.98"
$ws.Range("E48").Value = "  -1.25%  "

Set-TextCell $ws.Range("D49") "9.22"
$ws.Range("E49").Value = "  +6.39%  "

Set-TextCell $ws.Range("D50") "236.51"
$ws.Range("E50").Value = "  +15.83%  "

$ws.Range("E51").Value = "  -2.18%  "
